$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to track seat_min/seat_max per room; now it just records
# the last seat id used in a room (Seat.seat_last), so drop the old comment
# that warned about number formatting for the seat_min/seat_max columns.
$ws.Range("E1").Comment.Delete()

# Rename header E1 from seat_min -> seat_last
$ws.Range("E1").Value = "seat_last"

# Replace the old numeric seat_min values with the new textual seat_last ids
$ws.Range("E2").Value = "r2s7"
$ws.Range("E3").Value = "r2s3"

# Drop the seat_max column (F) entirely - no longer part of the schema
$ws.Range("F1:F3").Clear()

# Match the saved selection state
$null = $ws.Range("E4").Select()
